# Scheduled market-data refresh: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) on each job sheet from the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 636.1111
$ws.Range("J17").Value = 636.1111
$ws.Range("L17").Value = 1908.3333
$ws.Range("N17").Value = -2244.3333

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 4054.875
$ws.Range("I62").Value = 2111
$ws.Range("J62").Value = 4702.8335
$ws.Range("K62").Value = 2111
$ws.Range("L62").Value = 4702.8335
$ws.Range("M62").Value = -1487
$ws.Range("N62").Value = -5950.8335

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 3891.25
$ws.Range("I64").Value = 3279
$ws.Range("J64").Value = 4328.5713
$ws.Range("K64").Value = 3279
$ws.Range("L64").Value = 4328.5713
$ws.Range("M64").Value = -3031
$ws.Range("N64").Value = -4824.5713

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 4054.875
$ws.Range("I65").Value = 2111
$ws.Range("J65").Value = 4702.8335
$ws.Range("K65").Value = 10555
$ws.Range("L65").Value = 23514.1675
$ws.Range("M65").Value = -7435
$ws.Range("N65").Value = -29754.1675

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 3891.25
$ws.Range("I67").Value = 3279
$ws.Range("J67").Value = 4328.5713
$ws.Range("K67").Value = 3279
$ws.Range("L67").Value = 4328.5713
$ws.Range("M67").Value = -2421
$ws.Range("N67").Value = -6044.5713

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 1525.0526
$ws.Range("I86").Value = 1325.25
$ws.Range("J86").Value = 1670.3636
$ws.Range("K86").Value = 1325.25
$ws.Range("L86").Value = 1670.3636
$ws.Range("M86").Value = -202.25
$ws.Range("N86").Value = -3916.3636

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 1525.0526
$ws.Range("I89").Value = 1325.25
$ws.Range("J89").Value = 1670.3636
$ws.Range("K89").Value = 6626.25
$ws.Range("L89").Value = 8351.817999999999
$ws.Range("M89").Value = -1010.25
$ws.Range("N89").Value = -19583.818

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 1949.7222
$ws.Range("I116").Value = 1776.5385
$ws.Range("J116").Value = 2400
$ws.Range("K116").Value = 1776.5385
$ws.Range("L116").Value = 2400
$ws.Range("M116").Value = 1665.4615
$ws.Range("N116").Value = -9284

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1588557.8
$ws.Range("I137").Value = 2223199.8
$ws.Range("J137").Value = 1952.7778
$ws.Range("K137").Value = 6669599.399999999
$ws.Range("L137").Value = 5858.3334
$ws.Range("M137").Value = -6667049.399999999
$ws.Range("N137").Value = -10958.3334

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1813603.8
$ws.Range("I138").Value = 1469.8125
$ws.Range("J138").Value = 2780075.2
$ws.Range("K138").Value = 4409.4375
$ws.Range("L138").Value = 8340225.600000001
$ws.Range("M138").Value = 730.5625
$ws.Range("N138").Value = -8350505.600000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 1192.1
$ws.Range("I32").Value = 876.7619
$ws.Range("J32").Value = 2847.625
$ws.Range("K32").Value = 876.7619
$ws.Range("L32").Value = 2847.625
$ws.Range("M32").Value = -589.7619
$ws.Range("N32").Value = -3421.625

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 18556520
$ws.Range("I61").Value = 20429514
$ws.Range("J61").Value = 201182.8
$ws.Range("K61").Value = 20429514
$ws.Range("L61").Value = 201182.8
$ws.Range("M61").Value = -20429302
$ws.Range("N61").Value = -201606.8

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 6147741
$ws.Range("I74").Value = 8097961
$ws.Range("J74").Value = 102060
$ws.Range("K74").Value = 8097961
$ws.Range("L74").Value = 102060
$ws.Range("M74").Value = -8097087
$ws.Range("N74").Value = -103808

# Row 75: Someone Put Dung in My Helmet / Titanium Sallet of Maiming
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36748

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 6147741
$ws.Range("I77").Value = 8097961
$ws.Range("J77").Value = 102060
$ws.Range("K77").Value = 40489805
$ws.Range("L77").Value = 510300
$ws.Range("M77").Value = -40485437
$ws.Range("N77").Value = -519036

# Row 78: Rage against the Scream (L) / Titanium Sallet of Maiming
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -113736

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 46209.26
$ws.Range("I132").Value = 27072.842
$ws.Range("J132").Value = 137107.25
$ws.Range("K132").Value = 81218.526
$ws.Range("L132").Value = 411321.75
$ws.Range("M132").Value = -78688.526
$ws.Range("N132").Value = -416381.75

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 18556520
$ws.Range("I136").Value = 20429514
$ws.Range("J136").Value = 201182.8
$ws.Range("K136").Value = 61288542
$ws.Range("L136").Value = 603548.3999999999
$ws.Range("M136").Value = -61285992
$ws.Range("N136").Value = -608648.3999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 15: Anutha Spatha / Bronze Spatha
$ws.Range("H15").Value = 4000
$ws.Range("J15").Value = 5000
$ws.Range("L15").Value = 5000
$ws.Range("N15").Value = -5454

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1484.4717
$ws.Range("I134").Value = 898.23254
$ws.Range("J134").Value = 4005.3
$ws.Range("K134").Value = 2694.69762
$ws.Range("L134").Value = 12015.9
$ws.Range("M134").Value = -159.6976199999999
$ws.Range("N134").Value = -17085.9

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3331.5588
$ws.Range("I31").Value = 1206.1786
$ws.Range("K31").Value = 1206.1786
$ws.Range("M31").Value = -911.1786

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3331.5588
$ws.Range("I34").Value = 1206.1786
$ws.Range("K34").Value = 1206.1786
$ws.Range("M34").Value = -1004.1786

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 22751.94
$ws.Range("I134").Value = 1071.725
$ws.Range("K134").Value = 3215.175
$ws.Range("M134").Value = -680.1749999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 43: Sole Survivor / Baked Sole
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 115: Mixology / Blood Tomato Juice
$ws.Range("H115").Value = 2684.3333
$ws.Range("J115").Value = 2621.2
$ws.Range("L115").Value = 7863.599999999999
$ws.Range("N115").Value = -10213.6

# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 42886820
$ws.Range("I121").Value = 1189.8572
$ws.Range("J121").Value = 51225692
$ws.Range("K121").Value = 3569.5716
$ws.Range("L121").Value = 153677076
$ws.Range("M121").Value = -2259.5716
$ws.Range("N121").Value = -153679696

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1805.5883
$ws.Range("J132").Value = 1945
$ws.Range("L132").Value = 17505
$ws.Range("N132").Value = -22565

# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 5302.1113
$ws.Range("I133").Value = 2608.7778
$ws.Range("J133").Value = 7995.4443
$ws.Range("K133").Value = 7826.3334
$ws.Range("L133").Value = 23986.3329
$ws.Range("M133").Value = -2766.3334
$ws.Range("N133").Value = -34106.3329

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1030.5714
$ws.Range("I102").Value = 950
$ws.Range("J102").Value = 1138
$ws.Range("K102").Value = 950
$ws.Range("L102").Value = 1138
$ws.Range("M102").Value = 672
$ws.Range("N102").Value = -4382

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 625
$ws.Range("I122").Value = 625
$ws.Range("K122").Value = 1875
$ws.Range("M122").Value = 575

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 2035.3846
$ws.Range("I7").Value = 1795.7778
$ws.Range("K7").Value = 1795.7778
$ws.Range("M7").Value = -1683.7778

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 717.25
$ws.Range("I22").Value = 476.25
$ws.Range("J22").Value = 837.75
$ws.Range("K22").Value = 476.25
$ws.Range("L22").Value = 837.75
$ws.Range("M22").Value = -181.25
$ws.Range("N22").Value = -1427.75

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 717.25
$ws.Range("I27").Value = 476.25
$ws.Range("J27").Value = 837.75
$ws.Range("K27").Value = 476.25
$ws.Range("L27").Value = 837.75
$ws.Range("M27").Value = -369.25
$ws.Range("N27").Value = -1051.75

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 712.0513
$ws.Range("I46").Value = 636.6667
$ws.Range("J46").Value = 776.6667
$ws.Range("K46").Value = 636.6667
$ws.Range("L46").Value = 776.6667
$ws.Range("M46").Value = -448.6667
$ws.Range("N46").Value = -1152.6667

# Row 76: Dragoon Drop Rate / Dhalmelskin Breeches of Maiming
$ws.Range("H76").Value = 37788
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 37788
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 37788
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -38464

# Row 79: Exploiting the Adroit (L) / Dhalmelskin Breeches of Maiming
$ws.Range("H79").Value = 37788
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 37788
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 37788
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -40128

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3145.04
$ws.Range("I122").Value = 2271.75
$ws.Range("K122").Value = 6815.25
$ws.Range("M122").Value = -4365.25

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 2035.3846
$ws.Range("I126").Value = 1795.7778
$ws.Range("K126").Value = 5387.3334
$ws.Range("M126").Value = -2917.3334

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 43098.58
$ws.Range("I132").Value = 32201
$ws.Range("J132").Value = 66255.94
$ws.Range("K132").Value = 96603
$ws.Range("L132").Value = 198767.82
$ws.Range("M132").Value = -94073
$ws.Range("N132").Value = -203827.82

$ws = $wb.Worksheets.Item("WVR")
# Row 7: Long Hair, Long Life / Hempen Coif
$ws.Range("H7").Value = 2226
$ws.Range("I7").Value = 2226
$ws.Range("K7").Value = 2226
$ws.Range("M7").Value = -2113

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 2596.8147
$ws.Range("I122").Value = 2034
$ws.Range("J122").Value = 3300.3333
$ws.Range("K122").Value = 6102
$ws.Range("L122").Value = 9900.999899999999
$ws.Range("M122").Value = -3652
$ws.Range("N122").Value = -14800.9999

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 46267.6
$ws.Range("I132").Value = 37715.965
$ws.Range("J132").Value = 79252.5
$ws.Range("K132").Value = 113147.895
$ws.Range("L132").Value = 237757.5
$ws.Range("M132").Value = -110617.895
$ws.Range("N132").Value = -242817.5

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 81930.48
$ws.Range("I136").Value = 63391.375
$ws.Range("J136").Value = 114888.89
$ws.Range("K136").Value = 190174.125
$ws.Range("L136").Value = 344666.67
$ws.Range("M136").Value = -187624.125
$ws.Range("N136").Value = -349766.67
